$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 currently holds phone "09876543" stored as text (leading zero).
# Insert a new row 35 that is a duplicate of row 34 (phone kept as the
# original leading-zero text "09876543"), and convert row 34's phone
# value into a plain number (9876543, losing the leading zero).

# Shift row 34 down to make room, duplicating its contents into row 35.
$ws.Rows.Item(34).Copy()
$ws.Rows.Item(35).Insert()

# Row 34: phone becomes a true number.
$ws.Cells.Item(34, 1).Value = 9876543

# Row 35 was populated by the copy/insert above and already holds the
# original leading-zero phone text "09876543", a blank birthday, and a
# total_points of 0 - matching row 34's pre-edit contents, so nothing
# further is required there.
